$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.49%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '49.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.91%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.161'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.33%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07760'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.33%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.531'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.81%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.375'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '14.44%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.560'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-5.01%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1214'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-6.56%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1980'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.78%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04772'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.85%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09346'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.48%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.41%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001257'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-5.66%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005802'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.44%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,021.45%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.334'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.25%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.434'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.16%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3476'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.21%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.020'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.79%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1366'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.43%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.3038'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.77%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.43%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001270'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.59%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003927'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-7.60%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001351'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.15%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02602'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-4.85%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06088'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.44%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '76.50%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007900'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.73%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.30%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008388'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '8.91%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008350'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.06%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3378'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '5.76%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007534'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '7.45%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.03%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05319'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.09%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.03%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.03%'

Write-Host "Updated 71 cells (D/E columns) with new crypto price/volume data"